$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Native"
$ws.Range("I2:I98").Value = "Native"

$ws.Range("I2:I98").Select() | Out-Null
